$d = $word.ActiveDocument

# --- Edit 1: split the "Project Description" paragraph so that the
#     "The number of objects..." sentence starts its own paragraph.
$rng = $d.Content
$rng.Find.Execute("other particles that moves on the screen in horizontal directions.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertParagraphAfter()

# --- Edit 2: split the "Code URL" paragraph so the URL starts its own
#     paragraph, then collapse the run of 6 empty paragraphs down to 2.
$rng2 = $d.Content
$rng2.Find.Execute("serve as a useful resource for future development): ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0) | Out-Null
$rng2.InsertParagraphAfter()

$rng3 = $d.Content
$rng3.Find.Execute("https://github.com/zhxl0903/CSCB58-Project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$urlParaIndex = $rng3.Paragraphs.First.Index
for ($k = 0; $k -lt 4; $k++) {
    $emptyPara = $d.Paragraphs.Item($urlParaIndex + 1)
    $emptyPara.Range.Delete()
}

# --- Edit 3: fix capitalization in the "Week 1" update line.
$oldWeek1 = "Week 1: - Created repository (github for project work, google doc. for proposal) and organized method of communication (Facebook, google hangout)."
$newWeek1 = "Week 1: - Created repository (Github for project work, Google doc. for proposal) and organized method of communication (Facebook, google hangout)."
$d.Content.Find.Execute($oldWeek1, $true, $false, $false, $false, $false, $true, 1, $false, $newWeek1, 2) | Out-Null

# --- Edit 4: add three new update paragraphs after "Testing begins on the
#     restructured code."
$rng4 = $d.Content
$rng4.Find.Execute("Testing begins on the restructured code.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng4.Collapse(0) | Out-Null
$rng4.InsertParagraphAfter()
$testingIndex = $rng4.Paragraphs.First.Index

$p1 = $d.Paragraphs.Item($testingIndex + 1)
$p1.Range.InsertAfter("               - Game can now start with newly structured code. Graphic display is functional. Additional testing is required on collision detection and winning condition checking. ")

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($testingIndex + 2)
$p2.Range.InsertAfter("               - Documentations have been revised and errors have been corrected.")

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($testingIndex + 3)
$p3.Range.InsertAfter("             ")
